$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 23.51647972924649
$ws.Range("C2").Value = 16.07267116366805
$ws.Range("D2").Value = 6.09677355909968
$ws.Range("E2").Value = 12.32158857521591
$ws.Range("F2").Value = 48.68540974036421
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("J2").Value = 10.44424916278716
$ws.Range("N2").Value = 19.69065503691036
$ws.Range("B3").Value = 23.01450742996591
$ws.Range("C3").Value = 15.60482678393141
$ws.Range("D3").Value = 6.101291582118661
$ws.Range("E3").Value = 12.30726967054525
$ws.Range("F3").Value = 48.3830329837312
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("J3").Value = 10.45396233582873
$ws.Range("N3").Value = 19.75688350268325
$ws.Range("B4").Value = 22.7086730275475
$ws.Range("C4").Value = 15.31605601974711
$ws.Range("D4").Value = 6.104604796827756
$ws.Range("E4").Value = 12.30116094845516
$ws.Range("F4").Value = 48.21136809616186
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("J4").Value = 10.46192875252649
$ws.Range("N4").Value = 19.7995782131349
$ws.Range("B5").Value = 22.58483208960017
$ws.Range("C5").Value = 15.19820172182457
$ws.Range("D5").Value = 6.106091166166432
$ws.Range("E5").Value = 12.29934718394929
$ws.Range("F5").Value = 48.14497429615874
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("J5").Value = 10.46567755428517
$ws.Range("N5").Value = 19.81748732873598
$ws.Range("B6").Value = 22.56432167178713
$ws.Range("C6").Value = 15.17862747146456
$ws.Range("D6").Value = 6.106346221854106
$ws.Range("E6").Value = 12.29908683065347
$ws.Range("F6").Value = 48.134165823034
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("J6").Value = 10.46633034969084
$ws.Range("N6").Value = 19.82049197287164
$ws.Range("B7").Value = 22.70699941662838
$ws.Range("C7").Value = 15.31446704082637
$ws.Range("D7").Value = 6.104624290180083
$ws.Range("E7").Value = 12.30113375110292
$ws.Range("F7").Value = 48.21045821734276
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("J7").Value = 10.46197727745137
$ws.Range("N7").Value = 19.79981767375291
$ws.Range("B8").Value = 23.34301977724238
$ws.Range("C8").Value = 15.91179058953592
$ws.Range("D8").Value = 6.098219861105385
$ws.Range("E8").Value = 12.31609519766947
$ws.Range("F8").Value = 48.57827243357162
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("J8").Value = 10.44718201107515
$ws.Range("N8").Value = 19.71306935240545
$ws.Range("B9").Value = 24.60090203205951
$ws.Range("C9").Value = 17.06279902722094
$ws.Range("D9").Value = 6.089910612619341
$ws.Range("E9").Value = 12.36666690305405
$ws.Range("F9").Value = 49.40847033938102
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("J9").Value = 10.43410329698464
$ws.Range("N9").Value = 19.55905582560074
$ws.Range("B10").Value = 25.52091696912137
$ws.Range("C10").Value = 17.88574570553739
$ws.Range("D10").Value = 6.086359240176841
$ws.Range("E10").Value = 12.41666295741446
$ws.Range("F10").Value = 50.08155185634362
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("J10").Value = 10.43426526414429
$ws.Range("N10").Value = 19.45570354246046
$ws.Range("B11").Value = 25.93644470233297
$ws.Range("C11").Value = 18.25328407356964
$ws.Range("D11").Value = 6.085290033698149
$ws.Range("E11").Value = 12.44216427130292
$ws.Range("F11").Value = 50.4006549212662
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("J11").Value = 10.43646987167375
$ws.Range("N11").Value = 19.4108121546289
$ws.Range("B12").Value = 26.0932003413397
$ws.Range("C12").Value = 18.39133915044102
$ws.Range("D12").Value = 6.084963024510618
$ws.Range("E12").Value = 12.45221404385855
$ws.Range("F12").Value = 50.52327279289778
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("J12").Value = 10.43761158056307
$ws.Range("N12").Value = 19.39411844397326
$ws.Range("B13").Value = 26.05946923166854
$ws.Range("C13").Value = 18.36165862133882
$ws.Range("D13").Value = 6.085029999223311
$ws.Range("E13").Value = 12.45003222986603
$ws.Range("F13").Value = 50.49678691627506
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("J13").Value = 10.4373520394863
$ws.Range("N13").Value = 19.39770013340602
$ws.Range("B14").Value = 25.94935393400164
$ws.Range("C14").Value = 18.26466523656799
$ws.Range("D14").Value = 6.085261573653043
$ws.Range("E14").Value = 12.44298321850087
$ws.Range("F14").Value = 50.41070753745585
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("J14").Value = 10.43655764858525
$ws.Range("N14").Value = 19.4094326250779
$ws.Range("B15").Value = 25.88182281670292
$ws.Range("C15").Value = 18.20510365016561
$ws.Range("D15").Value = 6.085413540981623
$ws.Range("E15").Value = 12.43871655829139
$ws.Range("F15").Value = 50.35821095794817
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("J15").Value = 10.43611103555768
$ws.Range("N15").Value = 19.41665893213025
$ws.Range("B16").Value = 25.49368668051333
$ws.Range("C16").Value = 17.86157640403713
$ws.Range("D16").Value = 6.086440050457631
$ws.Range("E16").Value = 12.41505158465927
$ws.Range("F16").Value = 50.06095113938637
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("J16").Value = 10.43416410940235
$ws.Range("N16").Value = 19.45868006177299
$ws.Range("B17").Value = 25.25469061114688
$ws.Range("C17").Value = 17.64898006997492
$ws.Range("D17").Value = 6.087209223568896
$ws.Range("E17").Value = 12.40123782990431
$ws.Range("F17").Value = 49.8818464401668
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("J17").Value = 10.43351588324456
$ws.Range("N17").Value = 19.4850028485271
$ws.Range("B18").Value = 25.11695409256072
$ws.Range("C18").Value = 17.52606559685909
$ws.Range("D18").Value = 6.087703096274261
$ws.Range("E18").Value = 12.39355235122113
$ws.Range("F18").Value = 49.78005063558186
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("J18").Value = 10.4333435970522
$ws.Range("N18").Value = 19.50034293039593
$ws.Range("B19").Value = 25.07027731087183
$ws.Range("C19").Value = 17.48434442296574
$ws.Range("D19").Value = 6.08787917296984
$ws.Range("E19").Value = 12.39099490791819
$ws.Range("F19").Value = 49.74579621305429
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("J19").Value = 10.43331969342218
$ws.Range("N19").Value = 19.50557114685995
$ws.Range("B20").Value = 25.28016156330098
$ws.Range("C20").Value = 17.67167807489
$ws.Range("D20").Value = 6.087122022323816
$ws.Range("E20").Value = 12.40268146110948
$ws.Range("F20").Value = 49.90078666864979
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("J20").Value = 10.4335641274693
$ws.Range("N20").Value = 19.48218005131849
$ws.Range("B21").Value = 25.98171490233227
$ws.Range("C21").Value = 18.29318611028442
$ws.Range("D21").Value = 6.085191446593266
$ws.Range("E21").Value = 12.44504304963293
$ws.Range("F21").Value = 50.43594344694669
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("J21").Value = 10.4367826492582
$ws.Range("N21").Value = 19.40597820949631
$ws.Range("B22").Value = 26.43667881039389
$ws.Range("C22").Value = 18.69277101539929
$ws.Range("D22").Value = 6.084383302320571
$ws.Range("E22").Value = 12.47501738663038
$ws.Range("F22").Value = 50.79603903029474
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("J22").Value = 10.44067491858673
$ws.Range("N22").Value = 19.35795781012327
$ws.Range("B23").Value = 26.19423214126234
$ws.Range("C23").Value = 18.48015309726493
$ws.Range("D23").Value = 6.084773348114206
$ws.Range("E23").Value = 12.45881141894173
$ws.Range("F23").Value = 50.60292948250854
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("J23").Value = 10.43843376107117
$ws.Range("N23").Value = 19.38342406181026
$ws.Range("B24").Value = 25.26864717504147
$ws.Range("C24").Value = 17.66141844513771
$ws.Range("D24").Value = 6.087161285030811
$ws.Range("E24").Value = 12.40202799706612
$ws.Range("F24").Value = 49.89222012923249
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("J24").Value = 10.43354169212808
$ws.Range("N24").Value = 19.48345559393942
$ws.Range("B25").Value = 24.26059631667642
$ws.Range("C25").Value = 16.75471240777046
$ws.Range("D25").Value = 6.091707345534407
$ws.Range("E25").Value = 12.35072157389229
$ws.Range("F25").Value = 49.17253397481998
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("J25").Value = 10.43592908034525
$ws.Range("N25").Value = 19.59899853923204
